$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 30   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/20/2023  Through  11/26/2023"

# --- Cells that flip from a plain number to the shared "0" text ---
# (copy style+value from an existing s=14/t=s "0" cell, e.g. C14)
$ws.Range("C14").Copy($ws.Range("F14"))
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("C14").Copy($ws.Range("F28"))
$ws.Range("C14").Copy($ws.Range("F29"))
$ws.Range("C14").Copy($ws.Range("D23"))

# --- Cells that flip from a plain number to the shared "***.*" text ---
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("E14").Copy($ws.Range("E23"))

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -92.857142857142
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 150
$ws.Range("M15").Value = 31.578947368421
$ws.Range("N15").Value = -62.686567164179
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = -3.571428571428
$ws.Range("I16").Value = 341
$ws.Range("J16").Value = 360
$ws.Range("K16").Value = -5.277777777777
$ws.Range("L16").Value = 13.28903654485
$ws.Range("M16").Value = -19.953051643192
$ws.Range("N16").Value = -78.112965340179
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 59
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = 51.282051282051
$ws.Range("I17").Value = 604
$ws.Range("J17").Value = 529
$ws.Range("K17").Value = 14.177693761814
$ws.Range("L17").Value = 29.33618843683
$ws.Range("M17").Value = 89.937106918239
$ws.Range("N17").Value = -22.663252240717
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 149
$ws.Range("K18").Value = -22.797927461139
$ws.Range("L18").Value = -10.778443113772
$ws.Range("M18").Value = -34.361233480176
$ws.Range("N18").Value = -85.645472061657
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -12.5
$ws.Range("I19").Value = 534
$ws.Range("J19").Value = 547
$ws.Range("K19").Value = -2.376599634369
$ws.Range("L19").Value = 33.167082294264
$ws.Range("M19").Value = 41.269841269841
$ws.Range("N19").Value = -50.646950092421
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 208
$ws.Range("J20").Value = 198
$ws.Range("K20").Value = 5.050505050505
$ws.Range("L20").Value = 76.271186440678
$ws.Range("M20").Value = 35.064935064935
$ws.Range("N20").Value = -85.585585585585
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 158
$ws.Range("G21").Value = 138
$ws.Range("H21").Value = 14.492753623188
$ws.Range("I21").Value = 1863
$ws.Range("J21").Value = 1870
$ws.Range("K21").Value = -0.374331550802
$ws.Range("L21").Value = 25.963488843813
$ws.Range("M21").Value = 21.447196870925
$ws.Range("N21").Value = -68.934467233616
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 22
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 0
$ws.Range("L23").Value = -2.380952380952
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 123
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = 24.242424242424
$ws.Range("I24").Value = 1438
$ws.Range("J24").Value = 1346
$ws.Range("K24").Value = 6.835066864784
$ws.Range("L24").Value = 33.024976873265
$ws.Range("M24").Value = 40.566959921798
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 82
$ws.Range("G25").Value = 56
$ws.Range("H25").Value = 46.428571428571
$ws.Range("I25").Value = 900
$ws.Range("J25").Value = 708
$ws.Range("K25").Value = 27.118644067796
$ws.Range("L25").Value = 61.290322580645
$ws.Range("M25").Value = 26.050420168067
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 100
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -28.571428571428
$ws.Range("I27").Value = 81
$ws.Range("J27").Value = 76
$ws.Range("K27").Value = 6.578947368421
$ws.Range("L27").Value = 62
$ws.Range("H28").Value = -100
$ws.Range("M28").Value = -69.444444444444
$ws.Range("N28").Value = -92.617449664429
$ws.Range("H29").Value = -100
$ws.Range("M29").Value = -54.166666666666
$ws.Range("N29").Value = -91.970802919708
